# Updated Argent (Solar) prices: append a new row (row 47) with date
# 2025-04-17 and the latest price to every sheet in the workbook.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-17"

$values = @{
    "N-Dense"                    = "40"
    "N-Type"                     = "41"
    "N-type Wafer"                = "1.23"
    "Cell Topcon 183mm"          = "0.295"
    "Module Topcon 183mm"        = "0.09"
    "Silver Rear_side"            = "5,356"
    "Silver Busbar front-side"   = "8,018"
    "Silver finger front-side"   = "8,068"
    "USD_CNY"                    = "7.3413"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($values.ContainsKey($name)) {
        $cellA = $ws.Cells.Item(47, 1)
        $cellB = $ws.Cells.Item(47, 2)

        # Force text storage (so dates / numeric-looking strings are not
        # reinterpreted as a date serial or a number), then restore the
        # default "Normal" style so no extra number-format is left behind.
        $cellA.NumberFormat = "@"
        $cellB.NumberFormat = "@"
        $cellA.Value = $newDate
        $cellB.Value = $values[$name]
        $cellA.Style = "Normal"
        $cellB.Style = "Normal"
    }
}
